# Regenerate handback report timestamps (re-run of "Generate Report for Handback").
# The 729ef009-... row's timestamps are refreshed to reflect a new report run.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for 729ef009-...
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-02 10:53:44"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for 729ef009-...
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-02 10:53:40"
$wsZhCn.Range("K3").Value = "2016-09-02 10:53:58"

# de-de sheet: "Correspond Handoff Datetime" (shares text with Overview!G3) /
# "Correspond Handback DateTime" for 729ef009-...
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-02 10:53:44"
$wsDeDe.Range("K3").Value = "2016-09-02 10:54:14"
